$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header row: A1 = "relation", B1 = "count"
$ws.Range("A1").Value = "relation"
$ws.Range("B1").Value = "count"

# Widen column A to fit the new header text (closest width this host can
# represent to the authored 52.1640625 character-width value)
$ws.Columns.Item(1).ColumnWidth = 51.33

# Move selection to D5
$ws.Range("D5").Select()
